# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates metrics in row 3 of the sheet to reflect corrected relevance markers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.5729166666666666
$ws.Range("E3").Value = 0.9895833333333334
$ws.Range("H3").Value = 0.5574468085106383
$ws.Range("I3").Value = 0.1075262588816806
$ws.Range("J3").Value = 0.4791666666666667
$ws.Range("K3").Value = 164.5208333333333

$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 17
$ws.Range("S3").Value = 47
$ws.Range("T3").Value = 171
$ws.Range("U3").Value = 351
$ws.Range("V3").Value = 1073
$ws.Range("W3").Value = 1062
$ws.Range("X3").Value = 1032
$ws.Range("Y3").Value = 908
$ws.Range("Z3").Value = 728

$ws.Range("AF3").Value = 0.994439
$ws.Range("AG3").Value = 0.984245
$ws.Range("AH3").Value = 0.956441
$ws.Range("AI3").Value = 0.84152
$ws.Range("AJ3").Value = 0.674699
